# Added the New test cases
# Reproduces: B2 employee name change + a new data row (row 3) on the
# "Data" sheet of TestData.xlsx, matching the authors' commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Data" sheet

# ------------------------------------------------------------------
# 1) Existing row 2: Employee_Name changed from "Paul Collings" to "Mily Hm"
# ------------------------------------------------------------------
$ws.Range("B2").Value = "Mily Hm"

# ------------------------------------------------------------------
# 2) New row 3 with the new test case data
#    (set values first so the quote-prefixed date style carries over
#    correctly once formats are copied down below)
# ------------------------------------------------------------------
$ws.Range("A3").Value = "testcase2"
$ws.Range("B3").Value = "jabesh"
$ws.Range("C3").Value = "CAN - Personal"
$ws.Range("D3").Value = "2021-Apr-22"
$ws.Range("E3").Value = "2021-Apr-23"
$ws.Range("F3").Value = 2

# Copy the formatting (borders / alignment / quote-prefix) of row 2
# down onto row 3 so the new row matches the look of the existing data.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# TestCaseName (A3) and HowManyDayLeave (F3) keep the plain/default style.
$ws.Range("A3").ClearFormats()
$ws.Range("B3").ClearFormats()

# Employee_Name (B3) only needs horizontal centering (no border).
$ws.Range("B3").HorizontalAlignment = -4108   # xlHAlignCenter

# ------------------------------------------------------------------
# 3) Update selection to reflect the newly entered row
# ------------------------------------------------------------------
$ws.Range("A3:F3").Select()
